$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2450
$ws.Range("I40").Value = 2333.3333
$ws.Range("K40").Value = 2333.3333
$ws.Range("M40").Value = -2158.3333
$ws.Range("H64").Value = 3899.5
$ws.Range("I64").Value = 3899
$ws.Range("K64").Value = 3899
$ws.Range("M64").Value = -3651
$ws.Range("H67").Value = 3899.5
$ws.Range("I67").Value = 3899
$ws.Range("K67").Value = 3899
$ws.Range("M67").Value = -3041
$ws.Range("H70").Value = 7666.3335
$ws.Range("I70").Value = 4500
$ws.Range("J70").Value = 9249.5
$ws.Range("K70").Value = 13500
$ws.Range("L70").Value = 27748.5
$ws.Range("M70").Value = -13230
$ws.Range("N70").Value = -28288.5
$ws.Range("H73").Value = 7666.3335
$ws.Range("I73").Value = 4500
$ws.Range("J73").Value = 9249.5
$ws.Range("K73").Value = 13500
$ws.Range("L73").Value = 27748.5
$ws.Range("M73").Value = -12564
$ws.Range("N73").Value = -29620.5
$ws.Range("H98").Value = 1891.0698
$ws.Range("I98").Value = 1518.091
$ws.Range("K98").Value = 1518.091
$ws.Range("M98").Value = -20.09099999999989
$ws.Range("H107").Value = 880.7222
$ws.Range("I107").Value = 737.2308
$ws.Range("J107").Value = 1253.8
$ws.Range("K107").Value = 737.2308
$ws.Range("L107").Value = 1253.8
$ws.Range("M107").Value = 1182.7692
$ws.Range("N107").Value = -5093.8
$ws.Range("H113").Value = 19891.092
$ws.Range("I113").Value = 30043.572
$ws.Range("K113").Value = 30043.572
$ws.Range("M113").Value = -26789.572
$ws.Range("H122").Value = 1891.0698
$ws.Range("I122").Value = 1518.091
$ws.Range("K122").Value = 4554.272999999999
$ws.Range("M122").Value = -2104.272999999999
$ws.Range("H125").Value = 497
$ws.Range("I125").Value = 497
$ws.Range("K125").Value = 4473
$ws.Range("M125").Value = -2013
$ws.Range("H138").Value = 2418.6133
$ws.Range("J138").Value = 2282.7805
$ws.Range("L138").Value = 6848.3415
$ws.Range("N138").Value = -17128.3415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 198709.92
$ws.Range("I2").Value = 252784.36
$ws.Range("K2").Value = 252784.36
$ws.Range("M2").Value = -252671.36
$ws.Range("H5").Value = 193.33333
$ws.Range("I5").Value = 193.33333
$ws.Range("K5").Value = 193.33333
$ws.Range("M5").Value = -81.33332999999999
$ws.Range("H45").Value = 1381
$ws.Range("I45").Value = 1100.8182
$ws.Range("J45").Value = 1689.2
$ws.Range("K45").Value = 1100.8182
$ws.Range("L45").Value = 1689.2
$ws.Range("M45").Value = -723.8181999999999
$ws.Range("N45").Value = -2443.2
$ws.Range("H74").Value = 813.2222
$ws.Range("I74").Value = 521.619
$ws.Range("K74").Value = 521.619
$ws.Range("M74").Value = 352.381
$ws.Range("H77").Value = 813.2222
$ws.Range("I77").Value = 521.619
$ws.Range("K77").Value = 2608.095
$ws.Range("M77").Value = 1759.905
$ws.Range("H97").Value = 1795.3334
$ws.Range("I97").Value = 1793.3125
$ws.Range("J97").Value = 1801.8
$ws.Range("K97").Value = 1793.3125
$ws.Range("L97").Value = 1801.8
$ws.Range("M97").Value = -1297.3125
$ws.Range("N97").Value = -2793.8
$ws.Range("H116").Value = 198709.92
$ws.Range("I116").Value = 252784.36
$ws.Range("K116").Value = 252784.36
$ws.Range("M116").Value = -250490.36
$ws.Range("H123").Value = 81997.5
$ws.Range("J123").Value = 81997.5
$ws.Range("L123").Value = 81997.5
$ws.Range("N123").Value = -91797.5
$ws.Range("H132").Value = 2214.5557
$ws.Range("I132").Value = 1884.5652
$ws.Range("K132").Value = 5653.6956
$ws.Range("M132").Value = -3123.6956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 198709.92
$ws.Range("I3").Value = 252784.36
$ws.Range("K3").Value = 252784.36
$ws.Range("M3").Value = -252670.36
$ws.Range("H4").Value = 193.33333
$ws.Range("I4").Value = 193.33333
$ws.Range("K4").Value = 193.33333
$ws.Range("M4").Value = -78.33332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 883.8570999999999
$ws.Range("I16").Value = 837.4
$ws.Range("K16").Value = 837.4
$ws.Range("M16").Value = -550.4
$ws.Range("H62").Value = 2763.2
$ws.Range("I62").Value = 2668
$ws.Range("K62").Value = 2668
$ws.Range("M62").Value = -2044
$ws.Range("H65").Value = 2763.2
$ws.Range("I65").Value = 2668
$ws.Range("K65").Value = 13340
$ws.Range("M65").Value = -10220
$ws.Range("H113").Value = 883.8570999999999
$ws.Range("I113").Value = 837.4
$ws.Range("K113").Value = 837.4
$ws.Range("M113").Value = 1332.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 812.5
$ws.Range("J98").Value = 847.2222
$ws.Range("L98").Value = 2541.6666
$ws.Range("N98").Value = -5537.6666
$ws.Range("H118").Value = 1617.1428
$ws.Range("I118").Value = 1109.6666
$ws.Range("J118").Value = 1997.75
$ws.Range("K118").Value = 3328.9998
$ws.Range("L118").Value = 5993.25
$ws.Range("M118").Value = -2085.9998
$ws.Range("N118").Value = -8479.25
$ws.Range("H131").Value = 14762.775
$ws.Range("J131").Value = 15014.754
$ws.Range("L131").Value = 45044.262
$ws.Range("N131").Value = -55124.262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4378.1333
$ws.Range("I70").Value = 4079.375
$ws.Range("J70").Value = 4719.5713
$ws.Range("K70").Value = 4079.375
$ws.Range("L70").Value = 4719.5713
$ws.Range("M70").Value = -3809.375
$ws.Range("N70").Value = -5259.5713
$ws.Range("H73").Value = 4378.1333
$ws.Range("I73").Value = 4079.375
$ws.Range("J73").Value = 4719.5713
$ws.Range("K73").Value = 4079.375
$ws.Range("L73").Value = 4719.5713
$ws.Range("M73").Value = -3143.375
$ws.Range("N73").Value = -6591.5713
$ws.Range("H80").Value = 3374.875
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 2999.5
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 2999.5
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -4995.5
$ws.Range("H83").Value = 3374.875
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 2999.5
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 14997.5
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -24981.5
$ws.Range("H97").Value = 876.78125
$ws.Range("I97").Value = 874.76
$ws.Range("J97").Value = 884
$ws.Range("K97").Value = 874.76
$ws.Range("L97").Value = 884
$ws.Range("M97").Value = -378.76
$ws.Range("N97").Value = -1876
$ws.Range("H107").Value = 100
$ws.Range("I107").Value = 100
$ws.Range("K107").Value = 100
$ws.Range("M107").Value = 1820
$ws.Range("H122").Value = 1284.2
$ws.Range("I122").Value = 1239.5834
$ws.Range("J122").Value = 1351.125
$ws.Range("K122").Value = 3718.7502
$ws.Range("L122").Value = 4053.375
$ws.Range("M122").Value = -1268.7502
$ws.Range("N122").Value = -8953.375
$ws.Range("H132").Value = 1243305.4
$ws.Range("I132").Value = 1604885.8
$ws.Range("K132").Value = 4814657.4
$ws.Range("M132").Value = -4812127.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1980.9445
$ws.Range("I46").Value = 1555.7142
$ws.Range("K46").Value = 1555.7142
$ws.Range("M46").Value = -1367.7142
$ws.Range("H61").Value = 2136
$ws.Range("I61").Value = 1708.3334
$ws.Range("J61").Value = 2777.5
$ws.Range("K61").Value = 1708.3334
$ws.Range("L61").Value = 2777.5
$ws.Range("M61").Value = -1506.3334
$ws.Range("N61").Value = -3181.5
$ws.Range("H113").Value = 2136
$ws.Range("I113").Value = 1708.3334
$ws.Range("J113").Value = 2777.5
$ws.Range("K113").Value = 1708.3334
$ws.Range("L113").Value = 2777.5
$ws.Range("M113").Value = 461.6666
$ws.Range("N113").Value = -7117.5
$ws.Range("H132").Value = 2086.1956
$ws.Range("I132").Value = 1611.3478
$ws.Range("K132").Value = 4834.0434
$ws.Range("M132").Value = -2304.0434
$ws.Range("H136").Value = 2751.3684
$ws.Range("I136").Value = 2235.6875
$ws.Range("K136").Value = 6707.0625
$ws.Range("M136").Value = -4157.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 30892.8
$ws.Range("I16").Value = 42000
$ws.Range("K16").Value = 42000
$ws.Range("M16").Value = -41708
$ws.Range("H107").Value = 857
$ws.Range("I107").Value = 672.8570999999999
$ws.Range("K107").Value = 2018.5713
$ws.Range("M107").Value = -98.57129999999984
$ws.Range("H122").Value = 37991.184
$ws.Range("I122").Value = 58656.93
$ws.Range("K122").Value = 175970.79
$ws.Range("M122").Value = -173520.79
